$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Complete" shared string to "Completed" (cell C1)
$ws.Range("C1").Value = "Completed"

# Update column widths (ColumnWidth is offset by 5/6 from the stored xlsx
# "width" attribute in this engine, so compensate to land on the target)
$ws.Columns.Item(2).ColumnWidth = 23.98 - (5/6)
$ws.Columns.Item(3).ColumnWidth = 16.96 - (5/6)
$ws.Columns.Item(4).ColumnWidth = 16.22 - (5/6)
$ws.Columns.Item(5).ColumnWidth = 14.62 - (5/6)

# Update row heights for rows 2 and 3
$ws.Rows.Item(2).RowHeight = 25.1
$ws.Rows.Item(3).RowHeight = 25.1

# Update the active cell selection to E9
$ws.Range("E9").Select() | Out-Null
